{"js": "// Remove the sentence about sketching the t-distribution, and mark the\n// first row of the confidence-rating table as a repeating header row.\n\n// 1) Update the paragraph text that mentions the t-distribution applet.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText = \"Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.\";\nconst newText = \"Find the P-value and compare it to the level of significance.\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === oldText) {\n    p.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2) Mark the first row of the table as a header row (repeats on each page).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.headerRowCount = 1;\n}\nawait context.sync();\n", "ps1": "# Remove the sentence about sketching the t-distribution, and mark the\n# first row of the confidence-rating table as a repeating header row.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.\"\n$newText = \"Find the P-value and compare it to the level of significance.\"\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text.TrimEnd(\"`r\") -eq $oldText) {\n        $r.Text = $newText\n    }\n}\n\n$table = $d.Tables(1)\n$table.Rows(1).HeadingFormat = $true\n"}
